$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(51).Insert()

$ws.Range("A51").Value = 3
$ws.Range("B51").Value = "Femacal de La Calera"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 44536
$ws.Range("E51").Value = 5
$ws.Range("F51").Value = 100112010
$ws.Range("G51").Value = "Achicoria"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 78
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 5500
$ws.Range("M51").Value = 5256
$ws.Range("N51").Value = "$/caja 16 unidades"
$ws.Range("O51").Value = "Provincia de Quillota"
$ws.Range("P51").Value = 328
$ws.Range("Q51").Value = 16
$ws.Range("R51").Value = "Hortaliza"
